$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at AE (shifts AE:AM -> AF:AN)
$ws.Columns("AE:AE").Insert()

# New header cell value
$ws.Range("AE1").Value = "TIPO"

# Match new column's width to the "Cliente" (AD) column width
$ws.Columns("AE:AE").ColumnWidth = 25.5

# Select AE2 as the active cell (matches post-edit selection state)
$ws.Range("AE2").Select()
